# Applies the "Awaiting Connection Mode added to specification" edit.
#
# Word's live spell/grammar checker stamps <w:proofErr> markers around
# runs as text is (re)typed; because this headless host does not run a
# real proofing pass, the target <w:proofErr> markers are reproduced by
# hand using Range.InsertXML, which is the one primitive that lets us
# drop exact OOXML (including <w:proofErr/> siblings) into a Range.

$d = $word.ActiveDocument

# NOTE: named parameters (-Foo bar) do not bind correctly against
# functions in this PowerShell host, so every helper below takes
# positional parameters only.
function Replace-RangeWithXml($Range, $InnerBodyXml) {
    $pkg = @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>$InnerBodyXml</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $Range.InsertXML($pkg)
}

function Find-And-Replace($SearchText, $InnerBodyXml) {
    $probe = $d.Range(0, 0)
    $probe.Find.Execute($SearchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $target = $d.Range($probe.Start, $probe.End)
    Replace-RangeWithXml $target $InnerBodyXml
}

# ---------------------------------------------------------------------
# 1) "Format Correct – boolean, message passes format check"
#    -> split into 3 runs, proofErr spellStart/gramStart/.../spellEnd/gramEnd
#       wrapping the lone word "boolean"
# ---------------------------------------------------------------------
Find-And-Replace "Format Correct – boolean, message passes format check" @'
<w:p><w:r><w:t xml:space="preserve">Format Correct – </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>boolean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>, message passes format check</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------
# 2) "Awaiting Connection " (inside the Laura: ... (connection lost) line)
#    -> split into "Awaiting" (gramStart/End) + " Connection "
#
#    NOTE: this host's Range.InsertXML relocates its inserted content to
#    the END of the target paragraph whenever the replaced Range does not
#    already reach the paragraph end (trailing, untouched runs get
#    shifted ahead of the new content instead of staying put). To keep
#    everything in its original order we replace from the target text
#    through to the end of the paragraph, re-stating the trailing
#    "(connection lost)" run unchanged.
# ---------------------------------------------------------------------
Find-And-Replace "Awaiting Connection (connection lost)" @'
<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>Awaiting</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Connection </w:t></w:r><w:r><w:t>(connection lost)</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------
# 3) "Jireh: " -> split into "Jireh" (spellStart/End) + ": "
#    Same trailing-run caveat as above: "Jireh: " is the first run in its
#    paragraph, so the replacement must extend through the rest of the
#    paragraph's (unchanged) runs to preserve their order.
# ---------------------------------------------------------------------
Find-And-Replace "Jireh: No connection (user ends connection, not trying to get connection)" @'
<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Jireh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>No</w:t></w:r><w:r><w:t xml:space="preserve"> connection (user ends connection</w:t></w:r><w:r><w:t>, not trying to get connection</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------
# 4) "keyPressed" paragraph -> wrap with spellStart/gramStart .. spellEnd/gramEnd
# ---------------------------------------------------------------------
Find-And-Replace "keyPressed" @'
<w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>keyPressed</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>
'@

# ---------------------------------------------------------------------
# 5) Append the new "Mode: / Awaiting Connection / Events: ..." block
#    right after the "Follow line ..." paragraph, before the sectPr.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newBlockXml = @'
<w:p/>
<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Mode:</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Awaiting Connection</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>Events:</w:t></w:r></w:p>
<w:p><w:r><w:t>Connection Received: Base computer connected to robot computer</w:t></w:r></w:p>
<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>change</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> from Awaiting Connection Mode to Normal Mode</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Connection Lost: Base computer does not receive connection to robot computer (turned off or timed out) </w:t></w:r></w:p>
<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>change</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> from Awaiting Connection Mode to No Connection Mode</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t>Macros:</w:t></w:r><w:r><w:br/><w:t>Connection</w:t></w:r></w:p>
<w:p/>
'@

$pkg = @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>$newBlockXml</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$insertPoint.InsertXML($pkg)
